{"js": "// Replace the three-digit-by-one-digit multiplication results in the\n// table with newly generated problems/answers, as described in the\n// commit diff. Every \"old\" text value below is unique in the document,\n// so a plain search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"634\u00d77=4438\", \"306\u00d79=2754\"],\n  [\"251\u00d76=1506\", \"699\u00d72=1398\"],\n  [\"408\u00d74=1632\", \"499\u00d73=1497\"],\n  [\"524\u00d74=2096\", \"655\u00d74=2620\"],\n  [\"886\u00d72=1772\", \"879\u00d76=5274\"],\n  [\"763\u00d75=3815\", \"884\u00d74=3536\"],\n  [\"973\u00d74=3892\", \"135\u00d79=1215\"],\n  [\"529\u00d76=3174\", \"118\u00d73=354\"],\n  [\"494\u00d77=3458\", \"935\u00d74=3740\"],\n  [\"135\u00d78=1080\", \"656\u00d75=3280\"],\n  [\"711\u00d77=4977\", \"166\u00d78=1328\"],\n  [\"237\u00d75=1185\", \"545\u00d74=2180\"],\n  [\"730\u00d74=2920\", \"587\u00d73=1761\"],\n  [\"508\u00d79=4572\", \"346\u00d75=1730\"],\n  [\"117\u00d78=936\", \"750\u00d76=4500\"],\n  [\"748\u00d73=2244\", \"813\u00d77=5691\"],\n  [\"722\u00d75=3610\", \"859\u00d78=6872\"],\n  [\"243\u00d78=1944\", \"334\u00d79=3006\"],\n  [\"509\u00d79=4581\", \"907\u00d72=1814\"],\n  [\"902\u00d75=4510\", \"492\u00d79=4428\"],\n  [\"870\u00d74=3480\", \"860\u00d75=4300\"],\n  [\"715\u00d78=5720\", \"951\u00d75=4755\"],\n  [\"932\u00d76=5592\", \"206\u00d77=1442\"],\n  [\"921\u00d72=1842\", \"965\u00d76=5790\"],\n  [\"357\u00d79=3213\", \"886\u00d79=7974\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication results in the\n# table with newly generated problems/answers, as described in the\n# commit diff. Every \"old\" text value below is unique in the document,\n# so a plain Find/Replace per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"634\u00d77=4438\", \"306\u00d79=2754\"),\n    @(\"251\u00d76=1506\", \"699\u00d72=1398\"),\n    @(\"408\u00d74=1632\", \"499\u00d73=1497\"),\n    @(\"524\u00d74=2096\", \"655\u00d74=2620\"),\n    @(\"886\u00d72=1772\", \"879\u00d76=5274\"),\n    @(\"763\u00d75=3815\", \"884\u00d74=3536\"),\n    @(\"973\u00d74=3892\", \"135\u00d79=1215\"),\n    @(\"529\u00d76=3174\", \"118\u00d73=354\"),\n    @(\"494\u00d77=3458\", \"935\u00d74=3740\"),\n    @(\"135\u00d78=1080\", \"656\u00d75=3280\"),\n    @(\"711\u00d77=4977\", \"166\u00d78=1328\"),\n    @(\"237\u00d75=1185\", \"545\u00d74=2180\"),\n    @(\"730\u00d74=2920\", \"587\u00d73=1761\"),\n    @(\"508\u00d79=4572\", \"346\u00d75=1730\"),\n    @(\"117\u00d78=936\",  \"750\u00d76=4500\"),\n    @(\"748\u00d73=2244\", \"813\u00d77=5691\"),\n    @(\"722\u00d75=3610\", \"859\u00d78=6872\"),\n    @(\"243\u00d78=1944\", \"334\u00d79=3006\"),\n    @(\"509\u00d79=4581\", \"907\u00d72=1814\"),\n    @(\"902\u00d75=4510\", \"492\u00d79=4428\"),\n    @(\"870\u00d74=3480\", \"860\u00d75=4300\"),\n    @(\"715\u00d78=5720\", \"951\u00d75=4755\"),\n    @(\"932\u00d76=5592\", \"206\u00d77=1442\"),\n    @(\"921\u00d72=1842\", \"965\u00d76=5790\"),\n    @(\"357\u00d79=3213\", \"886\u00d79=7974\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
